# Update lemmatizer output: refreshed topic term lists and recalculated
# negativity/dominance scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Terms" column (D) text for each topic row.
$ws.Range("D2").Value = "payment,account,pay,month,check,make,interest,late,fee,bank,amount,charge,due,apply,send,take,time,balance,go,principal"
$ws.Range("D3").Value = "call,customer,service,number,speak,representative,phone,contact,information,give,account,issue,try,time,one,hold,supervisor,problem,person,someone"
$ws.Range("D4").Value = "tell,forbearance,receive,call,say,request,day,send,could,form,time,deferment,letter,state,ask,rep,month,payment,back,process"
$ws.Range("D5").Value = "credit,report,receive,send,account,letter,collection,agency,delinquent,state,contact,never,show,due,still,information,default,mail,issue,owe"
$ws.Range("D6").Value = "call,tell,say,payment,pay,time,day,phone,make,go,ask,harass,month,try,send,even,back,money,one,know"
$ws.Range("D7").Value = "pay,go,payment,year,make,month,take,job,work,time,interest,money,college,one,could,help,tell,try,like,even"

# Updated "Score for Negativity" (A) and "Percent Dominance" (B) values.
$ws.Range("A2").Value = -0.3288207330953067
$ws.Range("B2").Value = 0.1694656488549618

$ws.Range("A3").Value = -0.006713715523267633
$ws.Range("B3").Value = 0.07404580152671755

$ws.Range("A4").Value = 0.05709548479904948
$ws.Range("B4").Value = 0.07404580152671755

$ws.Range("A5").Value = -0.299271170347498
$ws.Range("B5").Value = 0.04274809160305344

$ws.Range("A6").Value = 0.4442888386644782
$ws.Range("B6").Value = 0.2572519083969466

$ws.Range("A7").Value = 0.1334213425059525
$ws.Range("B7").Value = 0.3824427480916031
